$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1131.9333
$ws.Range("I58").Value = 524.5
$ws.Range("J58").Value = 2346.8
$ws.Range("K58").Value = 1573.5
$ws.Range("L58").Value = 7040.400000000001
$ws.Range("M58").Value = -1423.5
$ws.Range("N58").Value = -7340.400000000001

$ws.Range("H107").Value = 2656.577
$ws.Range("I107").Value = 1716.9546
$ws.Range("K107").Value = 1716.9546
$ws.Range("M107").Value = 203.0454

$ws.Range("H132").Value = 5752934
$ws.Range("I132").Value = 6539683.5
$ws.Range("K132").Value = 19619050.5
$ws.Range("M132").Value = -19616520.5

$ws.Range("H138").Value = 2223.23
$ws.Range("J138").Value = 2384.7385
$ws.Range("L138").Value = 7154.2155
$ws.Range("N138").Value = -17434.2155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10166.0625
$ws.Range("I32").Value = 7490.0146
$ws.Range("J32").Value = 26952.182
$ws.Range("K32").Value = 7490.0146
$ws.Range("L32").Value = 26952.182
$ws.Range("M32").Value = -7203.0146
$ws.Range("N32").Value = -27526.182

$ws.Range("H61").Value = 55556760
$ws.Range("I61").Value = 83334264
$ws.Range("K61").Value = 83334264
$ws.Range("M61").Value = -83334052

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H122").Value = 3750.6667
$ws.Range("I122").Value = 3586.5334
$ws.Range("K122").Value = 10759.6002
$ws.Range("M122").Value = -8309.600199999999

$ws.Range("H132").Value = 2959.4524
$ws.Range("I132").Value = 2301.9
$ws.Range("J132").Value = 4603.3335
$ws.Range("K132").Value = 6905.700000000001
$ws.Range("L132").Value = 13810.0005
$ws.Range("M132").Value = -4375.700000000001
$ws.Range("N132").Value = -18870.0005

$ws.Range("H136").Value = 55556760
$ws.Range("I136").Value = 83334264
$ws.Range("K136").Value = 250002792
$ws.Range("M136").Value = -250000242

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 111123230
$ws.Range("I105").Value = 111123230
$ws.Range("K105").Value = 111123230
$ws.Range("M105").Value = -111121483

$ws.Range("H107").Value = 896.2593000000001
$ws.Range("I107").Value = 800.5909
$ws.Range("J107").Value = 1317.2
$ws.Range("K107").Value = 800.5909
$ws.Range("L107").Value = 1317.2
$ws.Range("M107").Value = 1119.4091
$ws.Range("N107").Value = -5157.2

$ws.Range("H134").Value = 9110.462
$ws.Range("I134").Value = 1610.875
$ws.Range("J134").Value = 21109.8
$ws.Range("K134").Value = 4832.625
$ws.Range("L134").Value = 63329.39999999999
$ws.Range("M134").Value = -2297.625
$ws.Range("N134").Value = -68399.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 50001084
$ws.Range("I16").Value = 66667764
$ws.Range("K16").Value = 66667764
$ws.Range("M16").Value = -66667477

$ws.Range("H28").Value = 19450
$ws.Range("J28").Value = 19450
$ws.Range("L28").Value = 19450
$ws.Range("N28").Value = -19940

$ws.Range("H31").Value = 1541.2632
$ws.Range("I31").Value = 1415.902
$ws.Range("J31").Value = 2606.8333
$ws.Range("K31").Value = 1415.902
$ws.Range("L31").Value = 2606.8333
$ws.Range("M31").Value = -1120.902
$ws.Range("N31").Value = -3196.8333

$ws.Range("H34").Value = 1541.2632
$ws.Range("I34").Value = 1415.902
$ws.Range("J34").Value = 2606.8333
$ws.Range("K34").Value = 1415.902
$ws.Range("L34").Value = 2606.8333
$ws.Range("M34").Value = -1213.902
$ws.Range("N34").Value = -3010.8333

$ws.Range("H113").Value = 50001084
$ws.Range("I113").Value = 66667764
$ws.Range("K113").Value = 66667764
$ws.Range("M113").Value = -66665594

$ws.Range("H141").Value = 684074.4399999999
$ws.Range("J141").Value = 684074.4399999999
$ws.Range("L141").Value = 684074.4399999999
$ws.Range("N141").Value = -694434.4399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 747.2308
$ws.Range("I68").Value = 1085
$ws.Range("J68").Value = 597.1111
$ws.Range("K68").Value = 3255
$ws.Range("L68").Value = 1791.3333
$ws.Range("M68").Value = -2444
$ws.Range("N68").Value = -3413.3333

$ws.Range("H71").Value = 747.2308
$ws.Range("I71").Value = 1085
$ws.Range("J71").Value = 597.1111
$ws.Range("K71").Value = 9765
$ws.Range("L71").Value = 5373.9999
$ws.Range("M71").Value = -5709
$ws.Range("N71").Value = -13485.9999

$ws.Range("H94").Value = 3803.4285
$ws.Range("I94").Value = 1412
$ws.Range("J94").Value = 4760
$ws.Range("K94").Value = 4236
$ws.Range("L94").Value = 14280
$ws.Range("M94").Value = -3560
$ws.Range("N94").Value = -15632

$ws.Range("H107").Value = 7218.2666
$ws.Range("J107").Value = 8866.333000000001
$ws.Range("L107").Value = 26598.999
$ws.Range("N107").Value = -30438.999

$ws.Range("H132").Value = 1206.0714
$ws.Range("I132").Value = 931.1111
$ws.Range("J132").Value = 1701
$ws.Range("K132").Value = 8379.999899999999
$ws.Range("L132").Value = 15309
$ws.Range("M132").Value = -5849.999899999999
$ws.Range("N132").Value = -20369

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 15200
$ws.Range("J64").Value = 15200
$ws.Range("L64").Value = 15200
$ws.Range("N64").Value = -15696

$ws.Range("H67").Value = 15200
$ws.Range("J67").Value = 15200
$ws.Range("L67").Value = 15200
$ws.Range("N67").Value = -16916

$ws.Range("H70").Value = 45003692
$ws.Range("I70").Value = 31254050
$ws.Range("J70").Value = 100002250
$ws.Range("K70").Value = 31254050
$ws.Range("L70").Value = 100002250
$ws.Range("M70").Value = -31253780
$ws.Range("N70").Value = -100002790

$ws.Range("H73").Value = 45003692
$ws.Range("I73").Value = 31254050
$ws.Range("J73").Value = 100002250
$ws.Range("K73").Value = 31254050
$ws.Range("L73").Value = 100002250
$ws.Range("M73").Value = -31253114
$ws.Range("N73").Value = -100004122

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 23176.666
$ws.Range("J64").Value = 23176.666
$ws.Range("L64").Value = 23176.666
$ws.Range("N64").Value = -23626.666

$ws.Range("H67").Value = 23176.666
$ws.Range("J67").Value = 23176.666
$ws.Range("L67").Value = 23176.666
$ws.Range("N67").Value = -24736.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2599.7
$ws.Range("I132").Value = 2178.5789
$ws.Range("J132").Value = 3327.0908
$ws.Range("K132").Value = 6535.736699999999
$ws.Range("L132").Value = 9981.2724
$ws.Range("M132").Value = -4005.736699999999
$ws.Range("N132").Value = -15041.2724
